# Auto-generated edit script: applies numeric cell updates described in the diff
# for Sheets/Pandaemonium_Profits.xlsx ("chore: update Sheets via scheduled runner")
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 733.3333
$ws.Range("I41").Value = 733.3333
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 733.3333
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -293.3333
$ws.Range("N41").ClearContents()
$ws.Range("H43").Value = 1399.9166
$ws.Range("I43").Value = 651
$ws.Range("J43").Value = 1468
$ws.Range("K43").Value = 651
$ws.Range("L43").Value = 1468
$ws.Range("M43").Value = -582
$ws.Range("N43").Value = -1606
$ws.Range("H53").Value = 222.78947
$ws.Range("I53").Value = 200.88889
$ws.Range("J53").Value = 242.5
$ws.Range("K53").Value = 200.88889
$ws.Range("L53").Value = 242.5
$ws.Range("M53").Value = 436.11111
$ws.Range("N53").Value = -1516.5
$ws.Range("H113").Value = 3339.25
$ws.Range("I113").Value = 2861.5
$ws.Range("K113").Value = 2861.5
$ws.Range("M113").Value = 392.5
$ws.Range("H138").Value = 2645.3943
$ws.Range("I138").Value = 1548.7567
$ws.Range("J138").Value = 3838.7942
$ws.Range("K138").Value = 4646.2701
$ws.Range("L138").Value = 11516.3826
$ws.Range("M138").Value = 493.7299000000003
$ws.Range("N138").Value = -21796.3826

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1488.225
$ws.Range("I45").Value = 1371
$ws.Range("J45").Value = 1892
$ws.Range("K45").Value = 1371
$ws.Range("L45").Value = 1892
$ws.Range("M45").Value = -994
$ws.Range("N45").Value = -2646
$ws.Range("H61").Value = 10068.154
$ws.Range("I61").Value = 9221.444
$ws.Range("J61").Value = 11973.25
$ws.Range("K61").Value = 9221.444
$ws.Range("L61").Value = 11973.25
$ws.Range("M61").Value = -9009.444
$ws.Range("N61").Value = -12397.25
$ws.Range("H132").Value = 3473.2812
$ws.Range("I132").Value = 1757.9286
$ws.Range("J132").Value = 4807.4443
$ws.Range("K132").Value = 5273.7858
$ws.Range("L132").Value = 14422.3329
$ws.Range("M132").Value = -2743.7858
$ws.Range("N132").Value = -19482.3329
$ws.Range("H136").Value = 10068.154
$ws.Range("I136").Value = 9221.444
$ws.Range("J136").Value = 11973.25
$ws.Range("K136").Value = 27664.332
$ws.Range("L136").Value = 35919.75
$ws.Range("M136").Value = -25114.332
$ws.Range("N136").Value = -41019.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 99.25
$ws.Range("I22").Value = 99.25
$ws.Range("K22").Value = 99.25
$ws.Range("M22").Value = 73.75
$ws.Range("H107").Value = 1641.8235
$ws.Range("I107").Value = 1344.4286
$ws.Range("J107").Value = 1850
$ws.Range("K107").Value = 1344.4286
$ws.Range("L107").Value = 1850
$ws.Range("M107").Value = 575.5714
$ws.Range("N107").Value = -5690
$ws.Range("H134").Value = 2627
$ws.Range("I134").Value = 1980.7778
$ws.Range("J134").Value = 3596.3333
$ws.Range("K134").Value = 5942.3334
$ws.Range("L134").Value = 10788.9999
$ws.Range("M134").Value = -3407.3334
$ws.Range("N134").Value = -15858.9999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3634.4348
$ws.Range("I122").Value = 3862.4614
$ws.Range("J122").Value = 3338
$ws.Range("K122").Value = 11587.3842
$ws.Range("L122").Value = 10014
$ws.Range("M122").Value = -9137.3842
$ws.Range("N122").Value = -14914

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 928.4286
$ws.Range("I7").Value = 933.1667
$ws.Range("J7").Value = 900
$ws.Range("K7").Value = 2799.5001
$ws.Range("L7").Value = 2700
$ws.Range("M7").Value = -2687.5001
$ws.Range("N7").Value = -2924
$ws.Range("H23").Value = 112.411766
$ws.Range("I23").Value = 90
$ws.Range("J23").Value = 113.8125
$ws.Range("K23").Value = 270
$ws.Range("L23").Value = 341.4375
$ws.Range("M23").Value = -35
$ws.Range("N23").Value = -811.4375
$ws.Range("H80").Value = 2490
$ws.Range("I80").Value = 2300
$ws.Range("J80").Value = 2571.4285
$ws.Range("K80").Value = 6900
$ws.Range("L80").Value = 7714.2855
$ws.Range("M80").Value = -5964
$ws.Range("N80").Value = -9586.2855
$ws.Range("H83").Value = 2490
$ws.Range("I83").Value = 2300
$ws.Range("J83").Value = 2571.4285
$ws.Range("K83").Value = 20700
$ws.Range("L83").Value = 23142.8565
$ws.Range("M83").Value = -16020
$ws.Range("N83").Value = -32502.8565
$ws.Range("H92").Value = 754.3333
$ws.Range("I92").Value = 667.3333
$ws.Range("J92").Value = 797.8333
$ws.Range("K92").Value = 2001.9999
$ws.Range("L92").Value = 2393.4999
$ws.Range("M92").Value = -753.9999
$ws.Range("N92").Value = -4889.4999
$ws.Range("H112").Value = 3258.4167
$ws.Range("I112").Value = 1084.6666
$ws.Range("J112").Value = 3983
$ws.Range("K112").Value = 3253.9998
$ws.Range("L112").Value = 11949
$ws.Range("M112").Value = -2145.9998
$ws.Range("N112").Value = -14165
$ws.Range("H113").Value = 675.13
$ws.Range("I113").Value = 689.50665
$ws.Range("J113").Value = 632
$ws.Range("K113").Value = 2068.51995
$ws.Range("L113").Value = 1896
$ws.Range("M113").Value = 101.4800500000001
$ws.Range("N113").Value = -6236
$ws.Range("H121").Value = 1083.3334
$ws.Range("I121").Value = 283
$ws.Range("J121").Value = 1483.5
$ws.Range("K121").Value = 849
$ws.Range("L121").Value = 4450.5
$ws.Range("M121").Value = 461
$ws.Range("N121").Value = -7070.5
$ws.Range("H131").Value = 1181.3334
$ws.Range("I131").Value = 851.6667
$ws.Range("J131").Value = 1304.9584
$ws.Range("K131").Value = 2555.0001
$ws.Range("L131").Value = 3914.8752
$ws.Range("M131").Value = 2484.9999
$ws.Range("N131").Value = -13994.8752

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 399.16666
$ws.Range("I22").Value = 433.33334
$ws.Range("J22").Value = 365
$ws.Range("K22").Value = 433.33334
$ws.Range("L22").Value = 365
$ws.Range("M22").Value = -138.33334
$ws.Range("N22").Value = -955
$ws.Range("H27").Value = 399.16666
$ws.Range("I27").Value = 433.33334
$ws.Range("J27").Value = 365
$ws.Range("K27").Value = 433.33334
$ws.Range("L27").Value = 365
$ws.Range("M27").Value = -326.33334
$ws.Range("N27").Value = -579
$ws.Range("H46").Value = 982.8333
$ws.Range("I46").Value = 1099.25
$ws.Range("J46").Value = 750
$ws.Range("K46").Value = 1099.25
$ws.Range("L46").Value = 750
$ws.Range("M46").Value = -911.25
$ws.Range("N46").Value = -1126

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 30000
$ws.Range("J97").Value = 30000
$ws.Range("L97").Value = 30000
$ws.Range("N97").Value = -31982

